# Update status text from "Ready for handoff" to "In Translation" on every
# sheet that shows a localization status (Overview summary columns + each
# per-locale sheet's Status column), then shrink the "Status"-column width
# back down (it had been widened to fit the old, longer text).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Text update -----------------------------------------------------------
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value     = "In Translation"
$wsDeDe.Range("C2").Value     = "In Translation"

# --- Column width update ----------------------------------------------------
# Target OOXML column width is 13.4101845877511 characters; Excel's
# ColumnWidth COM property is offset from the stored width by 5/6.
$newColumnWidth = 13.4101845877511 - (5/6)

$wsOverview.Range("E1:F1").ColumnWidth = $newColumnWidth
$wsZhCn.Range("C1").ColumnWidth        = $newColumnWidth
$wsDeDe.Range("C1").ColumnWidth        = $newColumnWidth
